$wb = $excel.ActiveWorkbook

# --- moving_average sheet ---
$ws = $wb.Worksheets.Item("moving_average")
$ws.Range("D14").Value = 24500
$ws.Range("D15").Value = 24500
$ws.Range("D16").Value = 24500
$ws.Range("D17").Value = 24500

# --- simple_exponential_smoothing sheet ---
$ws = $wb.Worksheets.Item("simple_exponential_smoothing")
$ws.Range("D15").Value = 23489.87524181193
$ws.Range("D16").Value = 23489.87524181193
$ws.Range("D17").Value = 23489.87524181193
$ws.Range("D18").Value = 23489.87524181193

# --- winter_trendseason sheet ---
$ws = $wb.Worksheets.Item("winter_trendseason")
$ws.Range("C3").Value = 18782.60227480593
$ws.Range("D3").Value = 505.7899893853545
$ws.Range("E3").Value = 0.4671052118983504
$ws.Range("C4").Value = 19264.45394236736
$ws.Range("D4").Value = 503.3961572029625
$ws.Range("E4").Value = 0.6825457909576081
$ws.Range("F4").Value = 13181.77283694694
$ws.Range("G4").Value = 181.7728369469369
$ws.Range("H4").Value = 181.7728369469369
$ws.Range("I4").Value = 462452.3122116908
$ws.Range("J4").Value = 563.0789798749702
$ws.Range("K4").Value = 1.398252591899515
$ws.Range("L4").Value = 6.601533313468529
$ws.Range("C5").Value = 19756.90420631653
$ws.Range("D5").Value = 502.3015678775828
$ws.Range("E5").Value = 1.1700523127902
$ws.Range("F5").Value = 23142.38273525652
$ws.Range("G5").Value = 142.3827352565168
$ws.Range("H5").Value = 142.3827352565168
$ws.Range("I5").Value = 315059.155907503
$ws.Range("J5").Value = 422.8468983354857
$ws.Range("K5").Value = 0.6190553706805076
$ws.Range("L5").Value = 4.607373999205856
$ws.Range("C6").Value = 20274.35573043596
$ws.Range("D6").Value = 503.8165635017684
$ws.Range("E6").Value = 1.66567736476266
$ws.Range("F6").Value = 33719.82347415544
$ws.Range("G6").Value = -280.1765258445594
$ws.Range("H6").Value = 280.1765258445594
$ws.Range("I6").Value = 255919.088339209
$ws.Range("J6").Value = 387.1793052127541
$ws.Range("K6").Value = 0.8240486054251746
$ws.Range("L6").Value = 3.661542650760685
$ws.Range("M6").Value = 2.552729848561492
$ws.Range("C7").Value = 20834.89756230393
$ws.Range("D7").Value = 509.4890903383878
$ws.Range("E7").Value = 0.4683910869955446
$ws.Range("F7").Value = 9705.592572220217
$ws.Range("G7").Value = -294.4074277797827
$ws.Range("H7").Value = 294.4074277797827
$ws.Range("I7").Value = 222070.4173777488
$ws.Range("J7").Value = 368.6249297261598
$ws.Range("K7").Value = 2.944074277797827
$ws.Range("L7").Value = 3.518048976168114
$ws.Range("M7").Value = 1.882555099834496
$ws.Range("C8").Value = 21796.85898455247
$ws.Range("D8").Value = 554.7363235294036
$ws.Range("E8").Value = 0.6968719177000449
$ws.Range("F8").Value = 14568.52127033276
$ws.Range("G8").Value = -3431.478729667237
$ws.Range("H8").Value = 3431.478729667237
$ws.Range("I8").Value = 2147566.393174569
$ws.Range("J8").Value = 879.1005630496726
$ws.Range("K8").Value = 19.06377072037354
$ws.Range("L8").Value = 6.109002600202351
$ws.Range("M8").Value = -3.114003224828376
$ws.Range("C9").Value = 22109.10339765365
$ws.Range("D9").Value = 530.4871324865809
$ws.Range("E9").Value = 1.157076628012107
$ws.Range("F9").Value = 26152.53578477179
$ws.Range("G9").Value = 3152.53578477179
$ws.Range("H9").Value = 3152.53578477179
$ws.Range("I9").Value = 3260554.319044871
$ws.Range("J9").Value = 1203.877023295689
$ws.Range("K9").Value = 13.70667732509474
$ws.Range("L9").Value = 7.194384703758407
$ws.Range("M9").Value = 0.3447310551293201
$ws.Range("C10").Value = 22655.24613435999
$ws.Range("D10").Value = 532.0526929085569
$ws.Range("E10").Value = 1.666841198160479
$ws.Range("F10").Value = 37710.25349354964
$ws.Range("G10").Value = -289.7465064503558
$ws.Range("H10").Value = 289.7465064503558
$ws.Range("I10").Value = 2863479.158914285
$ws.Range("J10").Value = 1089.610708690023
$ws.Range("K10").Value = 0.7624908064483048
$ws.Range("L10").Value = 6.390397966594644
$ws.Range("M10").Value = 0.1149651788820193
$ws.Range("C11").Value = 23406.20741380926
$ws.Range("D11").Value = 553.9435515626287
$ws.Range("E11").Value = 0.4728204293861066
$ws.Range("F11").Value = 10860.72410219483
$ws.Range("G11").Value = -1139.275897805168
$ws.Range("H11").Value = 1139.275897805168
$ws.Range("I11").Value = 2689531.426959339
$ws.Range("J11").Value = 1095.12906303615
$ws.Range("K11").Value = 9.493965815043067
$ws.Range("L11").Value = 6.735238838644468
$ws.Range("M11").Value = -0.9259261232257
$ws.Range("C12").Value = 23482.66857126982
$ws.Range("D12").Value = 506.1953121524218
$ws.Range("E12").Value = 0.6825447032705417
$ws.Range("F12").Value = 16697.15635162129
$ws.Range("G12").Value = 3697.156351621292
$ws.Range("H12").Value = 3697.156351621292
$ws.Range("I12").Value = 3787474.793096771
$ws.Range("J12").Value = 1355.331791894664
$ws.Range("K12").Value = 28.43966424324071
$ws.Range("L12").Value = 8.905681379104092
$ws.Range("M12").Value = 1.979698078285003
$ws.Range("C13").Value = 24318.89748263998
$ws.Range("D13").Value = 539.1986720741952
$ws.Range("E13").Value = 1.172953877819879
$ws.Range("F13").Value = 27756.95373207162
$ws.Range("G13").Value = -4243.046267928377
$ws.Range("H13").Value = 4243.046267928377
$ws.Range("I13").Value = 5079835.414795331
$ws.Range("J13").Value = 1617.851289715911
$ws.Range("K13").Value = 13.25951958727618
$ws.Range("L13").Value = 9.301484852574282
$ws.Range("M13").Value = -0.9641791764123474
$ws.Range("C14").Value = 24834.63567742981
$ws.Range("D14").Value = 536.8526243457586
$ws.Range("E14").Value = 1.665249091501127
$ws.Range("F14").Value = 41434.49877851217
$ws.Range("G14").Value = 434.4987785121702
$ws.Range("H14").Value = 434.4987785121702
$ws.Range("I14").Value = 4672248.2292731
$ws.Range("J14").Value = 1519.238580448932
$ws.Range("K14").Value = 1.059753118322366
$ws.Range("L14").Value = 8.614673874719955
$ws.Range("M14").Value = -0.7407656440841682
$ws.Range("E15").Value = 0.4728204293861066
$ws.Range("F15").Value = 11498.47154993854
$ws.Range("E16").Value = 0.6825447032705417
$ws.Range("F16").Value = 16966.76186378998
$ws.Range("E17").Value = 1.172953877819879
$ws.Range("F17").Value = 29789.85545321617
$ws.Range("E18").Value = 1.665249091501127
$ws.Range("F18").Value = 43190.72223570578
